$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the time-slot values in column C
$ws.Range("C6").Value = "20:10-20:15"
$ws.Range("C7").Value = "20:15-20:20"

# Move the active selection from C10 to C11
$ws.Range("C11").Select()
